$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C11").Value = "No"

$ws.Range("C7").Select()
